$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.231.62'
$ws.Range("E2").Value = '  +0.27%  '
$ws.Range("D3").Value = '1.834.23'
$ws.Range("E3").Value = '  -0.40%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9989'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.60'
$ws.Range("E5").Value = '  -0.67%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6243'
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("E7").Value = '  -0.21%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07381'
$ws.Range("E8").Value = '  -1.58%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2937'
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.21'
$ws.Range("E10").Value = '  -0.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07676'
$ws.Range("E11").Value = '  -0.52%  '
$ws.Range("D12").Value = '1.824.49'
$ws.Range("E12").Value = '  -3.19%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.967'
$ws.Range("E13").Value = '  -1.08%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6691'
$ws.Range("E14").Value = '  -1.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.66'
$ws.Range("E15").Value = '  -0.48%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008977'
$ws.Range("E16").Value = '  -3.35%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.882'
$ws.Range("E17").Value = '  -1.50%  '
$ws.Range("D18").Value = '29.195.61'
$ws.Range("E18").Value = '  +0.13%  '
$ws.Range("D19").Value = '2.075.09'
$ws.Range("E19").Value = '  -2.61%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '236.57'
$ws.Range("E20").Value = '  +2.62%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.49'
$ws.Range("E21").Value = '  -1.64%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.26%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.380'
$ws.Range("E23").Value = '  +2.65%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.001'
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.34'
$ws.Range("E25").Value = '  -1.44%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1409'
$ws.Range("E26").Value = '  +1.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.549'
$ws.Range("E27").Value = '  -0.10%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.68'
$ws.Range("E28").Value = '  -1.26%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.487'
$ws.Range("E29").Value = '  -1.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05857'
$ws.Range("E30").Value = '  +4.85%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.102'
$ws.Range("E31").Value = '  -1.39%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.087'
$ws.Range("E32").Value = '  -2.49%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.205'
$ws.Range("E33").Value = '  -0.09%  '
$ws.Range("E34").Value = '  +0.86%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7321'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.143'
$ws.Range("E36").Value = '  -0.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.608'
$ws.Range("E37").Value = '  -2.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.849'
$ws.Range("E38").Value = '  +2.74%  '
$ws.Range("D39").Value = '1.225.88'
$ws.Range("E39").Value = '  -0.03%  '
$ws.Range("E40").Value = '  -1.68%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.288'
$ws.Range("E41").Value = '  -4.31%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9143'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  -0.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.95'
$ws.Range("E44").Value = '  -0.26%  '
$ws.Range("D45").Value = '1.970.76'
$ws.Range("E45").Value = '  -2.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '65.19'
$ws.Range("E46").Value = '  -1.71%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5041'
$ws.Range("E47").Value = '  -1.22%  '
$ws.Range("E48").Value = '  -3.48%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4035'
$ws.Range("E49").Value = '  -1.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.115'
$ws.Range("E50").Value = '  -0.46%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1134'
$ws.Range("E51").Value = '  +2.88%  '
